$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record (Perejil, Vega Central Mapocho de Santiago) was
# added to the log. It belongs chronologically before the existing row 185,
# so insert a fresh row there (this pushes the old rows 185-205 down to
# 186-206, matching the diff) and fill it in with the new record's data.
$ws.Rows.Item(185).Insert()

$ws.Range("A185").Value = 9
$ws.Range("B185").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C185").Value = "Metropolitana"
$ws.Range("D185").Value = 44449
$ws.Range("E185").Value = 13
$ws.Range("F185").Value = 100112044
$ws.Range("G185").Value = "Perejil"
$ws.Range("H185").Value = "Sin especificar"
$ws.Range("I185").Value = "Primera"
$ws.Range("J185").Value = 106
$ws.Range("K185").Value = 9000
$ws.Range("L185").Value = 10000
$ws.Range("M185").Value = 9500
$ws.Range("N185").Value = '$/docena de atados'
$ws.Range("O185").Value = "Región Metropolitana"
$ws.Range("P185").Value = 3167
$ws.Range("Q185").Value = 3
$ws.Range("R185").Value = "Hortaliza"
